# Auto-generated edit script: applies the cryptos.xlsx price/volume update
# described by the commit "Updated cryptos list on Tue Dec 19 07:11:57 UTC 2023
# with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.959.13'
$ws.Range("E2").Value = '  +4.03%  '
$ws.Range("D3").Value = '2.246.11'
$ws.Range("E3").Value = '  +2.99%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'245.03"
$ws.Range("E5").Value = '  +2.89%  '
$ws.Range("D6").Value = "'0.618"
$ws.Range("E6").Value = '  +1.27%  '
$ws.Range("D7").Value = "'76.13"
$ws.Range("E7").Value = '  +8.58%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  +6.73%  '
$ws.Range("D10").Value = "'40.85"
$ws.Range("E10").Value = '  +1.75%  '
$ws.Range("D11").Value = "'0.0935"
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("D12").Value = "'6.95"
$ws.Range("E12").Value = '  +2.83%  '
$ws.Range("D13").Value = "'0.101"
$ws.Range("E13").Value = '  +0.02%  '
$ws.Range("D14").Value = '2.583.55'
$ws.Range("E14").Value = '  +3.23%  '
$ws.Range("E15").Value = '  +4.47%  '
$ws.Range("D16").Value = '2.232.87'
$ws.Range("E16").Value = '  +2.90%  '
$ws.Range("D17").Value = "'0.807"
$ws.Range("E17").Value = '  +0.91%  '
$ws.Range("D18").Value = '42.887.56'
$ws.Range("E18").Value = '  +4.30%  '
$ws.Range("E19").Value = '  +3.67%  '
$ws.Range("D20").Value = "'71.21"
$ws.Range("E20").Value = '  +0.90%  '
$ws.Range("E21").Value = '  +0.84%  '
$ws.Range("D22").Value = "'10.14"
$ws.Range("E22").Value = '  +3.59%  '
$ws.Range("D23").Value = "'231.00"
$ws.Range("E23").Value = '  +2.20%  '
$ws.Range("D24").Value = "'2.19"
$ws.Range("E24").Value = '  +13.18%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = "'10.87"
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").Value = "'3.43"
$ws.Range("E27").Value = '  -3.63%  '
$ws.Range("E28").Value = '  +2.18%  '
$ws.Range("D29").Value = "'38.59"
$ws.Range("E29").Value = '  +24.02%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = "'173.72"
$ws.Range("E30").Value = '  +3.45%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = "'2.14"
$ws.Range("E31").Value = '  -2.06%  '
$ws.Range("D32").Value = "'20.32"
$ws.Range("E32").Value = '  +1.79%  '
$ws.Range("D33").Value = "'0.0797"
$ws.Range("E33").Value = '  +3.53%  '
$ws.Range("D34").Value = "'5.34"
$ws.Range("E34").Value = '  +3.95%  '
$ws.Range("E35").Value = '  +1.15%  '
$ws.Range("D36").Value = "'0.109"
$ws.Range("E36").Value = '  +6.77%  '
$ws.Range("D37").Value = "'4.32"
$ws.Range("E37").Value = '  +4.65%  '
$ws.Range("D38").Value = "'0.0335"
$ws.Range("E38").Value = '  +17.01%  '
$ws.Range("E39").Value = '  +9.34%  '
$ws.Range("D40").Value = "'2.13"
$ws.Range("E40").Value = '  +2.27%  '
$ws.Range("D41").Value = "'5.55"
$ws.Range("E41").Value = '  +1.93%  '
$ws.Range("E42").Value = '  +6.49%  '
$ws.Range("D43").Value = "'106.27"
$ws.Range("E43").Value = '  +8.04%  '
$ws.Range("D44").Value = "'59.89"
$ws.Range("E44").Value = '  +0.43%  '
$ws.Range("D45").Value = "'8.67"
$ws.Range("E45").Value = '  +4.51%  '
$ws.Range("D46").Value = "'0.0992"
$ws.Range("E46").Value = '  +1.34%  '
$ws.Range("D47").Value = "'0.459"
$ws.Range("E47").Value = '  +23.09%  '
$ws.Range("E48").Value = '  +5.80%  '
$ws.Range("E49").Value = '  +1.55%  '
$ws.Range("E50").Value = '  +2.00%  '
$ws.Range("D51").Value = '2.450.68'
$ws.Range("E51").Value = '  +2.89%  '
